$d = $word.ActiveDocument

# Paragraph 2 is the (currently empty) paragraph directly below
# "URL to GitHub Repository:" that should receive the repo URL. Its
# paragraph-mark run properties are presently <w:b/><w:szCs w:val="24"/>
# and need to become <w:bCs/><w:szCs w:val="24"/> to match the rest of
# the "answer" paragraphs in this document (e.g. the YouTube URL
# paragraph right below). The Word object model has no direct setter
# for the paragraph-mark's "bCs" flag, so instead we clone a paragraph
# that already carries the desired mark formatting (the YouTube URL
# paragraph), fill it with the GitHub URL text, move it into place, and
# discard the old empty paragraph.

# 1) Clone the paragraph-mark formatting of the YouTube URL paragraph
#    (currently paragraph 4) by inserting a brand-new empty paragraph
#    right after it; the new paragraph inherits its <w:bCs/><w:szCs/>
#    mark formatting cleanly.
$videoUrlPara = $d.Paragraphs.Item(4)
$videoUrlPara.Range.InsertParagraphAfter() | Out-Null

# 2) Fill the new paragraph (now index 5) with the GitHub repo URL.
$newPara = $d.Paragraphs.Item(5)
$newPara.Range.InsertAfter("https://github.com/Cmendence/week_6_completed.git")

# 3) Cut the whole new paragraph (text + its own mark) to the clipboard.
$newPara2 = $d.Paragraphs.Item(5)
$newPara2.Range.Cut()

# 4) Paste it back in right before the old empty paragraph (index 2).
$target = $d.Paragraphs.Item(2)
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)
$insertionPoint.Paste()

# 5) Delete the now-redundant old empty paragraph (shifted to index 3).
$old = $d.Paragraphs.Item(3)
$old.Range.Delete()

Write-Output "done"
